$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("info")

# Add new config row (dbdriver / JDBC driver class)
$ws.Range("A8").Value = "dbdriver"
$ws.Range("B8").Value = "com.microsoft.sqlserver.jdbc.SQLServerDriver"

# Widen column B to fit the new longer value; columns C:D keep their width
$ws.Columns.Item(2).ColumnWidth = 42 - (5/7)

# Update the active selection on the sheet
$ws.Range("B3").Select()
